# Updated symbol list refresh: pushes new Price (column D) and
# Volume(1h) (column E) quotes into the Sheet1 table for rows 2-51.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the sheet's existing inline-string cells)
# instead of re-parsing look-alike numbers/percentages into numeric
# values, then the style is restored to Normal so the quote-prefix
# flag doesn't leave a stray cell format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '309.20'
    'E2' = '0.48%'
    'D3' = '41.31'
    'E3' = '3.48%'
    'D4' = '5.126'
    'E4' = '0.71%'
    'D5' = '0.07653'
    'E5' = '-0.43%'
    'D6' = '4.267'
    'E6' = '0.56%'
    'D7' = '1.620'
    'E7' = '-0.12%'
    'D8' = '2.473'
    'E8' = '2.13%'
    'D9' = '0.9085'
    'E9' = '-0.95%'
    'D10' = '0.1164'
    'E10' = '12.16%'
    'D11' = '0.1805'
    'E11' = '1.49%'
    'D12' = '0.09239'
    'E12' = '-0.64%'
    'D13' = '0.04249'
    'E13' = '-4.26%'
    'D14' = '0.1042'
    'E14' = '-1.17%'
    'E15' = '-0.94%'
    'D16' = '0.005822'
    'E16' = '-0.43%'
    'D17' = '3.356'
    'E17' = '-0.17%'
    'E18' = '0.61%'
    'D19' = '6.910'
    'E19' = '0.04%'
    'D20' = '0.1368'
    'E20' = '1.66%'
    'D21' = '0.2735'
    'E21' = '1.12%'
    'D22' = '0.04061'
    'E22' = '-2.20%'
    'E23' = '5.89%'
    'D24' = '0.004067'
    'E24' = '-1.12%'
    'D25' = '0.0001270'
    'E25' = '-2.23%'
    'D26' = '0.0003744'
    'D38' = '0.02433'
    'E38' = '-2.22%'
    'D39' = '0.05269'
    'E39' = '1.39%'
    'D40' = '0.007782'
    'E40' = '-1.86%'
    'D41' = '0.1305'
    'E41' = '-0.94%'
    'D42' = '0.006785'
    'E42' = '-4.11%'
    'E43' = '0.07%'
    'D44' = '0.007565'
    'E44' = '2.02%'
    'D45' = '0.3084'
    'E45' = '0.46%'
    'D46' = '0.00006891'
    'E46' = '7.57%'
    'E47' = '0.08%'
    'D48' = '0.07743'
    'E48' = '1,607.96%'
    'D49' = '0.003000'
    'E49' = '0.03%'
    'E50' = '0.08%'
    'E51' = '0.08%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}

Write-Output "Updated $($updates.Count) cells on $($ws.Name)."
